# Moderation decision dates that were "April 30, 2025" have now rolled over
# to "May 1, 2025". Column G ("Moderation Decision Date") is the only
# column affected; update every matching cell in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # column G
    if ($cell.Value2 -eq "April 30, 2025") {
        $cell.Value = "May 1, 2025"
    }
}
